$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.863899999999999
$ws.Range("C4").Value = -12.09099999999999
$ws.Range("B6").Value = 5.484899999999999
$ws.Range("B7").Value = 5.216500000000003
$ws.Range("D7").Value = -8.194499999999998
$ws.Range("D8").Value = -8.7546
$ws.Range("C9").Value = -10.1406
$ws.Range("D10").Value = -7.712100000000001
$ws.Range("C12").Value = -10.54389999999999
$ws.Range("D13").Value = -8.252299999999996
$ws.Range("B16").Value = 4.576399999999999
$ws.Range("D16").Value = -8.313499999999994
$ws.Range("C17").Value = -13.88389999999999
$ws.Range("C18").Value = -11.7807
$ws.Range("C19").Value = -11.0409
$ws.Range("B20").Value = 9.76259999999999
$ws.Range("C20").Value = -12.67399999999999
$ws.Range("C26").Value = -12.5857
$ws.Range("B28").Value = 6.333300000000003
$ws.Range("B29").Value = 5.5127
$ws.Range("D30").Value = -6.918699999999996
$ws.Range("C31").Value = -13.2754
$ws.Range("B32").Value = 6.573899999999997
$ws.Range("C39").Value = -11.5255
$ws.Range("B40").Value = 9.430399999999993
$ws.Range("C40").Value = -12.65560000000001
$ws.Range("D40").Value = -8.655199999999995
$ws.Range("C41").Value = -12.60510000000001
$ws.Range("C42").Value = -11.68729999999999
$ws.Range("C43").Value = -12.40899999999999
$ws.Range("D44").Value = -6.518300000000004
$ws.Range("B46").Value = 5.293700000000004
$ws.Range("C47").Value = -12.19629999999999
$ws.Range("C48").Value = -11.5375
$ws.Range("B51").Value = 5.329599999999997
$ws.Range("B52").Value = 5.257499999999994
$ws.Range("B57").Value = 5.289599999999998
$ws.Range("B59").Value = 4.741200000000003
$ws.Range("B62").Value = 5.197199999999999
$ws.Range("C63").Value = -10.1728
$ws.Range("C64").Value = -10.3897
$ws.Range("B66").Value = 5.775799999999999
$ws.Range("B73").Value = 9.021799999999995
$ws.Range("B74").Value = 9.194999999999991
$ws.Range("C76").Value = -11.7392
$ws.Range("C81").Value = -13.53879999999999
$ws.Range("C89").Value = -12.746
$ws.Range("D89").Value = -7.804899999999996
$ws.Range("D91").Value = -7.770199999999997
$ws.Range("B92").Value = 4.805699999999997
$ws.Range("C94").Value = -10.6441
$ws.Range("B100").Value = 4.681900000000003
